$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "English"
$ws.Cells.Item(1, 2).Value = "Japanese"
$ws.Cells.Item(2, 1).Value = "(my) older brother"
$ws.Cells.Item(2, 2).Value = "兄|あに"
$ws.Cells.Item(3, 1).Value = "landlord; landlady"
$ws.Cells.Item(3, 2).Value = "大家さん|おおやさん"
$ws.Cells.Item(4, 1).Value = "return (as a token of gratitude)"
$ws.Cells.Item(4, 2).Value = "お返し|おかえし"
$ws.Cells.Item(5, 1).Value = "(your/his) wife"
$ws.Cells.Item(5, 2).Value = "奥さん|おくさん"
$ws.Cells.Item(6, 1).Value = "uncle; middle-aged man"
$ws.Cells.Item(6, 2).Value = "おじさん"
$ws.Cells.Item(7, 1).Value = "aunt; middle-aged woman"
$ws.Cells.Item(7, 2).Value = "おばさん"
$ws.Cells.Item(8, 1).Value = "tumbler; glass"
$ws.Cells.Item(8, 2).Value = "グラス"
$ws.Cells.Item(9, 1).Value = "Christmas"
$ws.Cells.Item(9, 2).Value = "クリスマス"
$ws.Cells.Item(10, 1).Value = "(your/his) husband"
$ws.Cells.Item(10, 2).Value = "ご主人|ごしゅじん"
$ws.Cells.Item(11, 1).Value = "plate; dish"
$ws.Cells.Item(11, 2).Value = "皿|さら"
$ws.Cells.Item(12, 1).Value = "time"
$ws.Cells.Item(12, 2).Value = "時間|じかん"
$ws.Cells.Item(13, 1).Value = "ticket"
$ws.Cells.Item(13, 2).Value = "チケット"
$ws.Cells.Item(14, 1).Value = "chocolate"
$ws.Cells.Item(14, 2).Value = "チョコレート"
$ws.Cells.Item(15, 1).Value = "sweat shirt"
$ws.Cells.Item(15, 2).Value = "トレーナー"
$ws.Cells.Item(16, 1).Value = "stuffed animal"
$ws.Cells.Item(16, 2).Value = "ぬいぐるみ"
$ws.Cells.Item(17, 1).Value = "necktie"
$ws.Cells.Item(17, 2).Value = "ネクタイ"
$ws.Cells.Item(18, 1).Value = "St. Valentine's Day"
$ws.Cells.Item(18, 2).Value = "バレンタインデー"
$ws.Cells.Item(19, 1).Value = "camcorder"
$ws.Cells.Item(19, 2).Value = "ビデオカメラ"
$ws.Cells.Item(20, 1).Value = "married couple; husband and wife"
$ws.Cells.Item(20, 2).Value = "夫婦|ふうふ"
$ws.Cells.Item(21, 1).Value = "White Day"
$ws.Cells.Item(21, 2).Value = "ホワイトデー"
$ws.Cells.Item(22, 1).Value = "winter scarf"
$ws.Cells.Item(22, 2).Value = "マフラー"
$ws.Cells.Item(23, 1).Value = "comic book"
$ws.Cells.Item(23, 2).Value = "漫画|まんが"
$ws.Cells.Item(24, 1).Value = "multistory apartment building; condo"
$ws.Cells.Item(24, 2).Value = "マンション"
$ws.Cells.Item(25, 1).Value = "mandarin orange"
$ws.Cells.Item(25, 2).Value = "みかん"
$ws.Cells.Item(26, 1).Value = "everyone; all of you"
$ws.Cells.Item(26, 2).Value = "皆さん|みなさん"
$ws.Cells.Item(27, 1).Value = "ring"
$ws.Cells.Item(27, 2).Value = "指輪|ゆびわ"
$ws.Cells.Item(28, 1).Value = "radio"
$ws.Cells.Item(28, 2).Value = "ラジオ"
$ws.Cells.Item(29, 1).Value = "parents"
$ws.Cells.Item(29, 2).Value = "両親|りょうしん"
$ws.Cells.Item(30, 1).Value = "résumé"
$ws.Cells.Item(30, 2).Value = "履歴書|りれきしょ"
$ws.Cells.Item(31, 1).Value = "to want"
$ws.Cells.Item(31, 2).Value = "欲しい|ほしい"
$ws.Cells.Item(32, 1).Value = "stingy; cheap"
$ws.Cells.Item(32, 2).Value = "けち（な）"
$ws.Cells.Item(33, 1).Value = "to send"
$ws.Cells.Item(33, 2).Value = "送る|おくる"
$ws.Cells.Item(34, 1).Value = "to look good (on somebody)"
$ws.Cells.Item(34, 2).Value = "似合う|にあう"
$ws.Cells.Item(35, 1).Value = "to give up"
$ws.Cells.Item(35, 2).Value = "あきらめる"
$ws.Cells.Item(36, 1).Value = "to give (to others)"
$ws.Cells.Item(36, 2).Value = "あげる"
$ws.Cells.Item(37, 1).Value = "to give (me)"
$ws.Cells.Item(37, 2).Value = "くれる"
$ws.Cells.Item(38, 1).Value = "to come into existence; to be made"
$ws.Cells.Item(38, 2).Value = "できる"
$ws.Cells.Item(39, 1).Value = "to consult"
$ws.Cells.Item(39, 2).Value = "相談する|そうだんする"
$ws.Cells.Item(40, 1).Value = "to propose marriage"
$ws.Cells.Item(40, 2).Value = "プロポーズする"
$ws.Cells.Item(41, 1).Value = "same"
$ws.Cells.Item(41, 2).Value = "同じ|おなじ"
$ws.Cells.Item(42, 1).Value = "Mr./Ms....(casual)"
$ws.Cells.Item(42, 2).Value = "～君|～くん"
$ws.Cells.Item(43, 1).Value = "...like this; this kind of..."
$ws.Cells.Item(43, 2).Value = "こんな～"
$ws.Cells.Item(44, 1).Value = "[makes a noun plural]"
$ws.Cells.Item(44, 2).Value = "～たち"
$ws.Cells.Item(45, 1).Value = "we"
$ws.Cells.Item(45, 2).Value = "私たち|わたしたち"
$ws.Cells.Item(46, 1).Value = "exactly"
$ws.Cells.Item(46, 2).Value = "ちょうど"
$ws.Cells.Item(47, 1).Value = "what should one do"
$ws.Cells.Item(47, 2).Value = "どうしたらいい"
$ws.Cells.Item(48, 1).Value = "well"
$ws.Cells.Item(48, 2).Value = "よく"
$ws.Cells.Item(49, 1).Value = "[generic counter for smaller items]"
$ws.Cells.Item(49, 2).Value = "～個|～こ"
$ws.Cells.Item(50, 1).Value = "[counter for bound volumes]"
$ws.Cells.Item(50, 2).Value = "～冊|～さつ"
$ws.Cells.Item(51, 1).Value = "[counter for equipment]"
$ws.Cells.Item(51, 2).Value = "～台|～だい"
$ws.Cells.Item(52, 1).Value = "[counter for smaller animals]"
$ws.Cells.Item(52, 2).Value = "～匹|～ひき"
$ws.Cells.Item(53, 1).Value = "[counter for long objects]"
$ws.Cells.Item(53, 2).Value = "～本|～ほん"
$ws.Cells.Item(54, 1).Value = "[counter for flat objects]"
$ws.Cells.Item(54, 2).Value = "～枚|～まい"
$ws.Cells.Item(55, 1).Value = "he; boyfriend"
$ws.Cells.Item(55, 2).Value = "彼|かれ"
$ws.Cells.Item(56, 1).Value = "she; girlfriend"
$ws.Cells.Item(56, 2).Value = "彼女|かのじょ"
$ws.Cells.Item(57, 1).Value = "they"
$ws.Cells.Item(57, 2).Value = "彼ら|かれら"
$ws.Cells.Item(58, 1).Value = "age; era"
$ws.Cells.Item(58, 2).Value = "時代|じだい"
$ws.Cells.Item(59, 1).Value = "electricity fee"
$ws.Cells.Item(59, 2).Value = "電気代|でんきだい"
$ws.Cells.Item(60, 1).Value = "90's"
$ws.Cells.Item(60, 2).Value = "九十年代|きゅうじゅうねんだい"
$ws.Cells.Item(61, 1).Value = "in one's teens"
$ws.Cells.Item(61, 2).Value = "十代|じゅうだい"
$ws.Cells.Item(62, 1).Value = "instead"
$ws.Cells.Item(62, 2).Value = "代わりに|かわりに"
$ws.Cells.Item(63, 1).Value = "foreign students"
$ws.Cells.Item(63, 2).Value = "留学生|りゅうがくせい"
$ws.Cells.Item(64, 1).Value = "to study abroad"
$ws.Cells.Item(64, 2).Value = "留学する|りゅうがくする"
$ws.Cells.Item(65, 1).Value = "absence from home"
$ws.Cells.Item(65, 2).Value = "留守|るす"
$ws.Cells.Item(66, 1).Value = "family"
$ws.Cells.Item(66, 2).Value = "家族|かぞく"
$ws.Cells.Item(67, 1).Value = "race"
$ws.Cells.Item(67, 2).Value = "民族|みんぞく"
$ws.Cells.Item(68, 1).Value = "aquarium"
$ws.Cells.Item(68, 2).Value = "水族館|すいぞくかん"
$ws.Cells.Item(69, 1).Value = "member of royalty"
$ws.Cells.Item(69, 2).Value = "王族|おうぞく"
$ws.Cells.Item(70, 1).Value = "father"
$ws.Cells.Item(70, 2).Value = "父親|ちちおや"
$ws.Cells.Item(71, 1).Value = "kind"
$ws.Cells.Item(71, 2).Value = "親切な|しんせつな"
$ws.Cells.Item(72, 1).Value = "best friend"
$ws.Cells.Item(72, 2).Value = "親友|しんゆう"
$ws.Cells.Item(73, 1).Value = "parents"
$ws.Cells.Item(73, 2).Value = "両親|りょうしん"
$ws.Cells.Item(74, 1).Value = "intimate"
$ws.Cells.Item(74, 2).Value = "親しい|したしい"
$ws.Cells.Item(75, 1).Value = "mother"
$ws.Cells.Item(75, 2).Value = "母親|ははおや"
$ws.Cells.Item(76, 1).Value = "to cut"
$ws.Cells.Item(76, 2).Value = "切る|きる"
$ws.Cells.Item(77, 1).Value = "ticket"
$ws.Cells.Item(77, 2).Value = "切符|きっぷ"
$ws.Cells.Item(78, 1).Value = "postage stamp"
$ws.Cells.Item(78, 2).Value = "切手|きって"
$ws.Cells.Item(79, 1).Value = "precious"
$ws.Cells.Item(79, 2).Value = "大切な|たいせつな"
$ws.Cells.Item(80, 1).Value = "English"
$ws.Cells.Item(80, 2).Value = "英語|えいご"
$ws.Cells.Item(81, 1).Value = "England"
$ws.Cells.Item(81, 2).Value = "英国|えいこく"
$ws.Cells.Item(82, 1).Value = "English conversation"
$ws.Cells.Item(82, 2).Value = "英会話|えいかいわ"
$ws.Cells.Item(83, 1).Value = "hero"
$ws.Cells.Item(83, 2).Value = "英雄|えいゆう"
$ws.Cells.Item(84, 1).Value = "shop"
$ws.Cells.Item(84, 2).Value = "店|みせ"
$ws.Cells.Item(85, 1).Value = "store attendant"
$ws.Cells.Item(85, 2).Value = "店員|てんいん"
$ws.Cells.Item(86, 1).Value = "stall"
$ws.Cells.Item(86, 2).Value = "売店|ばいてん"
$ws.Cells.Item(87, 1).Value = "book store"
$ws.Cells.Item(87, 2).Value = "書店|しょてん"
$ws.Cells.Item(88, 1).Value = "last year"
$ws.Cells.Item(88, 2).Value = "去年|きょねん"
$ws.Cells.Item(89, 1).Value = "the past"
$ws.Cells.Item(89, 2).Value = "過去|かこ"
$ws.Cells.Item(90, 1).Value = "to leave"
$ws.Cells.Item(90, 2).Value = "去る|さる"
$ws.Cells.Item(91, 1).Value = "to erase"
$ws.Cells.Item(91, 2).Value = "消去する|しょうきょする"
$ws.Cells.Item(92, 1).Value = "suddenly"
$ws.Cells.Item(92, 2).Value = "急に|きゅうに"
$ws.Cells.Item(93, 1).Value = "to hurry"
$ws.Cells.Item(93, 2).Value = "急ぐ|いそぐ"
$ws.Cells.Item(94, 1).Value = "express train"
$ws.Cells.Item(94, 2).Value = "急行|きゅうこう"
$ws.Cells.Item(95, 1).Value = "super express"
$ws.Cells.Item(95, 2).Value = "特急|とっきゅう"
$ws.Cells.Item(96, 1).Value = "to ride"
$ws.Cells.Item(96, 2).Value = "乗る|のる"
$ws.Cells.Item(97, 1).Value = "vehicle"
$ws.Cells.Item(97, 2).Value = "乗り物|のりもの"
$ws.Cells.Item(98, 1).Value = "riding a car"
$ws.Cells.Item(98, 2).Value = "乗車|じょうしゃ"
$ws.Cells.Item(99, 1).Value = "horseback riding"
$ws.Cells.Item(99, 2).Value = "乗馬|じょうば"
$ws.Cells.Item(100, 1).Value = "really"
$ws.Cells.Item(100, 2).Value = "本当に|ほんとうに"
$ws.Cells.Item(101, 1).Value = "lunch box"
$ws.Cells.Item(101, 2).Value = "お弁当|おべんとう"
$ws.Cells.Item(102, 1).Value = "at that time"
$ws.Cells.Item(102, 2).Value = "当時|とうじ"
$ws.Cells.Item(103, 1).Value = "to hit"
$ws.Cells.Item(103, 2).Value = "当たる|あたる"
$ws.Cells.Item(104, 1).Value = "music"
$ws.Cells.Item(104, 2).Value = "音楽|おんがく"
$ws.Cells.Item(105, 1).Value = "pronunciation"
$ws.Cells.Item(105, 2).Value = "発音|はつおん"
$ws.Cells.Item(106, 1).Value = "sound"
$ws.Cells.Item(106, 2).Value = "音|おと"
$ws.Cells.Item(107, 1).Value = "real intention"
$ws.Cells.Item(107, 2).Value = "本音|ほんね"
$ws.Cells.Item(108, 1).Value = "fun"
$ws.Cells.Item(108, 2).Value = "楽しい|たのしい"
$ws.Cells.Item(109, 1).Value = "musical instrument"
$ws.Cells.Item(109, 2).Value = "楽器|がっき"
$ws.Cells.Item(110, 1).Value = "easy; comfortable"
$ws.Cells.Item(110, 2).Value = "楽な|らくな"
$ws.Cells.Item(111, 1).Value = "doctor"
$ws.Cells.Item(111, 2).Value = "医者|いしゃ"
$ws.Cells.Item(112, 1).Value = "dentist"
$ws.Cells.Item(112, 2).Value = "歯医者|はいしゃ"
$ws.Cells.Item(113, 1).Value = "medical science"
$ws.Cells.Item(113, 2).Value = "医学|いがく"
$ws.Cells.Item(114, 1).Value = "clinic"
$ws.Cells.Item(114, 2).Value = "医院|いいん"
$ws.Cells.Item(115, 1).Value = "scholar"
$ws.Cells.Item(115, 2).Value = "学者|がくしゃ"
$ws.Cells.Item(116, 1).Value = "reader"
$ws.Cells.Item(116, 2).Value = "読者|どくしゃ"
$ws.Cells.Item(117, 1).Value = "young people"
$ws.Cells.Item(117, 2).Value = "若者|わかもの"
